$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / "want to go" counts)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14859
$ws1.Range("F3").Value = 18453
$ws1.Range("F17").Value = 1406
$ws1.Range("F22").Value = 7643
$ws1.Range("F25").Value = 51
$ws1.Range("F26").Value = 1215
$ws1.Range("F28").Value = 5948
$ws1.Range("F29").Value = 97
$ws1.Range("F31").Value = 156
$ws1.Range("F34").Value = 5285

# Sheet "全部类型" (All types) - same events, shifted down by 1 row after row 19
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14859
$ws4.Range("F3").Value = 18453
$ws4.Range("F17").Value = 1406
$ws4.Range("F23").Value = 7643
$ws4.Range("F26").Value = 51
$ws4.Range("F27").Value = 1215
$ws4.Range("F30").Value = 5948
$ws4.Range("F31").Value = 97
$ws4.Range("F33").Value = 156
$ws4.Range("F36").Value = 5285
